{"js": "// Update the date heading and the 25 division-problem answers in the table.\n//\n// Each answer cell is addressed by its (row, column) position in the single\n// body table rather than by searching for its old text, because one answer\n// (\"60\u00f74=15, 0\") appears twice in the original table with two different\n// replacements \u2014 a text search could not tell those two occurrences apart.\n// The cell's existing paragraph is rewritten with\n// Paragraph.insertText(text, \"Replace\") (not by clearing the cell body),\n// which keeps that paragraph's formatting (font, size, alignment) exactly\n// as it was \u2014 only the literal <w:t> text changes, matching the diff.\n\nconst body = context.document.body;\n\n// 1) Date heading paragraph (first paragraph in the body, above the table).\nconst headingParas = body.paragraphs;\nheadingParas.load(\"items\");\nawait context.sync();\nheadingParas.items[0].insertText(\"2025-10-14 Tuesday\", \"Replace\");\n\n// 2) Table cell contents \u2014 keyed by 0-based (row, col) in the body's table.\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\nconst newValues = {\n  0: [\"23\u00f74=5, 3\", \"90\u00f73=30, 0\", \"59\u00f75=11, 4\", \"64\u00f72=32, 0\", \"82\u00f72=41, 0\"],\n  4: [\"84\u00f75=16, 4\", \"32\u00f75=6, 2\", \"64\u00f77=9, 1\", \"74\u00f78=9, 2\", \"36\u00f78=4, 4\"],\n  8: [\"95\u00f76=15, 5\", \"40\u00f75=8, 0\", \"21\u00f74=5, 1\", \"22\u00f76=3, 4\", \"25\u00f73=8, 1\"],\n  12: [\"40\u00f79=4, 4\", \"14\u00f77=2, 0\", \"98\u00f77=14, 0\", \"16\u00f75=3, 1\", \"35\u00f78=4, 3\"],\n  16: [\"17\u00f72=8, 1\", \"39\u00f79=4, 3\", \"58\u00f74=14, 2\", \"90\u00f72=45, 0\", \"94\u00f78=11, 6\"],\n};\n\n// Collect every cell's first paragraph up front, then load them all in one\n// sync round-trip before rewriting any text.\nconst cellParas = [];\nfor (const rowIndex of Object.keys(newValues)) {\n  const r = Number(rowIndex);\n  const vals = newValues[r];\n  for (let c = 0; c < vals.length; c++) {\n    const para = table.getCell(r, c).body.paragraphs.getFirst();\n    para.load(\"text\");\n    cellParas.push({ r, c, para, text: vals[c] });\n  }\n}\nawait context.sync();\n\nfor (const { para, text } of cellParas) {\n  para.insertText(text, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Update the date heading and the 25 division-problem answers in the table.\n# Cells are addressed by (row, column) rather than by searching for their old\n# text, because one answer (\"60\u00f74=15, 0\") appears twice in the original table\n# with two different replacements \u2014 a text-based Find/Replace would not be\n# able to tell those two occurrences apart. Assigning Cell.Range.Text keeps\n# each cell paragraph's existing formatting (font, size, alignment) intact,\n# matching the diff (only the <w:t> text itself changes).\n\n$d = $word.ActiveDocument\n\n# 1) Date heading, above the table.\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"2025-10-13 Monday\"\n$find.Replacement.Text = \"2025-10-14 Tuesday\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n# 2) Table answers \u2014 keyed by 1-based (row, column) in the single body table.\n$table = $d.Tables.Item(1)\n\n$newValues = @{\n    1  = @(\"23\u00f74=5, 3\", \"90\u00f73=30, 0\", \"59\u00f75=11, 4\", \"64\u00f72=32, 0\", \"82\u00f72=41, 0\")\n    5  = @(\"84\u00f75=16, 4\", \"32\u00f75=6, 2\", \"64\u00f77=9, 1\", \"74\u00f78=9, 2\", \"36\u00f78=4, 4\")\n    9  = @(\"95\u00f76=15, 5\", \"40\u00f75=8, 0\", \"21\u00f74=5, 1\", \"22\u00f76=3, 4\", \"25\u00f73=8, 1\")\n    13 = @(\"40\u00f79=4, 4\", \"14\u00f77=2, 0\", \"98\u00f77=14, 0\", \"16\u00f75=3, 1\", \"35\u00f78=4, 3\")\n    17 = @(\"17\u00f72=8, 1\", \"39\u00f79=4, 3\", \"58\u00f74=14, 2\", \"90\u00f72=45, 0\", \"94\u00f78=11, 6\")\n}\n\nforeach ($row in $newValues.Keys) {\n    $vals = $newValues[$row]\n    for ($col = 1; $col -le $vals.Length; $col++) {\n        $table.Cell($row, $col).Range.Text = $vals[$col - 1]\n    }\n}\n"}
